$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.965.30"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "2.910.95"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").Value = "'572.58"
$ws.Range("E5").Value = "  -3.90%  "

$ws.Range("D6").Value = "'144.77"
$ws.Range("E6").Value = "  -1.30%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.503"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "2.904.44"
$ws.Range("E9").Value = "  -2.07%  "

$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "  -6.96%  "

$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").Value = "'0.434"
$ws.Range("E12").Value = "  -2.76%  "

$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "  -2.49%  "

$ws.Range("D14").Value = "'32.31"
$ws.Range("E14").Value = "  -3.19%  "

$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("D16").Value = "3.372.13"
$ws.Range("E16").Value = "  -2.59%  "

$ws.Range("D17").Value = "61.755.26"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").Value = "'6.64"
$ws.Range("E18").Value = "  -1.60%  "

$ws.Range("D19").Value = "2.909.08"
$ws.Range("E19").Value = "  -2.33%  "

$ws.Range("D20").Value = "'437.15"
$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").Value = "'13.29"
$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").Value = "'0.658"
$ws.Range("E22").Value = "  -2.39%  "

$ws.Range("D23").Value = "'6.93"
$ws.Range("E23").Value = "  -2.39%  "

$ws.Range("D24").Value = "'79.47"
$ws.Range("E24").Value = "  -2.73%  "

$ws.Range("D25").Value = "'12.01"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").Value = "'10.23"
$ws.Range("E26").Value = "  -9.30%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'2.05"
$ws.Range("E28").Value = "  -4.60%  "

$ws.Range("D29").Value = "'0.0000110"
$ws.Range("E29").Value = "  +12.94%  "

$ws.Range("D30").Value = "'7.07"
$ws.Range("E30").Value = "  -2.77%  "

$ws.Range("D31").Value = "'2.52"
$ws.Range("E31").Value = "  -3.94%  "

$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  -4.15%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("D35").Value = "'25.72"
$ws.Range("E35").Value = "  -3.39%  "

$ws.Range("D36").Value = "'0.966"
$ws.Range("E36").Value = "  -2.75%  "

$ws.Range("D37").Value = "'5.46"
$ws.Range("E37").Value = "  -3.64%  "

$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  -4.35%  "

$ws.Range("D39").Value = "'49.10"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").Value = "'1.97"
$ws.Range("E40").Value = "  -4.02%  "

$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("D42").Value = "'8.29"
$ws.Range("E42").Value = "  -2.89%  "

$ws.Range("D43").Value = "'0.272"
$ws.Range("E43").Value = "  -4.11%  "

$ws.Range("D44").Value = "'38.87"
$ws.Range("E44").Value = "  -5.45%  "

$ws.Range("D45").Value = "2.687.43"
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("D46").Value = "'133.50"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").Value = "'0.0335"
$ws.Range("E47").Value = "  -1.77%  "

$ws.Range("D49").Value = "'339.92"
$ws.Range("E49").Value = "  -7.38%  "

$ws.Range("D50").Value = "'0.103"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").Value = "'21.79"
$ws.Range("E51").Value = "  -5.55%  "
